# Natmi following Dr Hou advice
# Update recomputed NATMI LR-pair stats (Tnfsf14-Tnfrsf14) for rows 2-16,
# columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T with the re-run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 0.2853273333333333
$ws.Cells.Item(2, 8).Value = 0.855982
$ws.Cells.Item(2, 9).Value = 0.05293626506635817
$ws.Cells.Item(2, 10).Value = 0.05293626506635819
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.461605
$ws.Cells.Item(2, 14).Value = 10.384815
$ws.Cells.Item(2, 15).Value = 0.06616740110392129
$ws.Cells.Item(2, 16).Value = 0.06733956725325306
$ws.Cells.Item(2, 17).Value = 0.9876905237033333
$ws.Cells.Item(2, 18).Value = 8.88921471333
$ws.Cells.Item(2, 19).Value = 0.003502655083589217
$ws.Cells.Item(2, 20).Value = 0.003564705181572058

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 0.2853273333333333
$ws.Cells.Item(3, 8).Value = 0.855982
$ws.Cells.Item(3, 9).Value = 0.05293626506635817
$ws.Cells.Item(3, 10).Value = 0.05293626506635819
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.998462333333334
$ws.Cells.Item(3, 14).Value = 11.995387
$ws.Cells.Item(3, 15).Value = 0.07642924626252497
$ws.Cells.Item(3, 16).Value = 0.07778320264879995
$ws.Cells.Item(3, 17).Value = 1.140870595003778
$ws.Cells.Item(3, 18).Value = 10.267835355034
$ws.Cells.Item(3, 19).Value = 0.004045878838974986
$ws.Cells.Item(3, 20).Value = 0.004117552233127129

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 0.2853273333333333
$ws.Cells.Item(4, 8).Value = 0.855982
$ws.Cells.Item(4, 9).Value = 0.05293626506635817
$ws.Cells.Item(4, 10).Value = 0.05293626506635819
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 21.74728766666667
$ws.Cells.Item(4, 14).Value = 65.241863
$ws.Cells.Item(4, 15).Value = 0.4156920000874432
$ws.Cells.Item(4, 16).Value = 0.4230560507063459
$ws.Cells.Item(4, 17).Value = 6.205095597162888
$ws.Cells.Item(4, 18).Value = 55.845860374466
$ws.Cells.Item(4, 19).Value = 0.02200518190259348
$ws.Cells.Item(4, 20).Value = 0.0223950072381178

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 0.2853273333333333
$ws.Cells.Item(5, 8).Value = 0.855982
$ws.Cells.Item(5, 9).Value = 0.05293626506635817
$ws.Cells.Item(5, 10).Value = 0.05293626506635819
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 20.37655766666667
$ws.Cells.Item(5, 14).Value = 61.129673
$ws.Cells.Item(5, 15).Value = 0.3894909627896643
$ws.Cells.Item(5, 16).Value = 0.3963908578200832
$ws.Cells.Item(5, 17).Value = 5.813988861542889
$ws.Cells.Item(5, 18).Value = 52.325899753886
$ws.Cells.Item(5, 19).Value = 0.02061819684718472
$ws.Cells.Item(5, 20).Value = 0.02098345151944503

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 0.2853273333333333
$ws.Cells.Item(6, 8).Value = 0.855982
$ws.Cells.Item(6, 9).Value = 0.05293626506635817
$ws.Cells.Item(6, 10).Value = 0.05293626506635819
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 2.731955
$ws.Cells.Item(6, 14).Value = 5.46391
$ws.Cells.Item(6, 15).Value = 0.0522203897564463
$ws.Cells.Item(6, 16).Value = 0.03543032157151783
$ws.Cells.Item(6, 17).Value = 0.7795014349366667
$ws.Cells.Item(6, 18).Value = 4.677008609620001
$ws.Cells.Item(6, 19).Value = 0.002764352394015776
$ws.Cells.Item(6, 20).Value = 0.001875548894096176

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.792877
$ws.Cells.Item(7, 8).Value = 5.378630999999999
$ws.Cells.Item(7, 9).Value = 0.3326292332200106
$ws.Cells.Item(7, 10).Value = 0.3326292332200106
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.461605
$ws.Cells.Item(7, 14).Value = 10.384815
$ws.Cells.Item(7, 15).Value = 0.06616740110392129
$ws.Cells.Item(7, 16).Value = 0.06733956725325306
$ws.Cells.Item(7, 17).Value = 6.206231987584999
$ws.Cells.Item(7, 18).Value = 55.85608788826499
$ws.Cells.Item(7, 19).Value = 0.02200921189335822
$ws.Cells.Item(7, 20).Value = 0.0223991086208169

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.792877
$ws.Cells.Item(8, 8).Value = 5.378630999999999
$ws.Cells.Item(8, 9).Value = 0.3326292332200106
$ws.Cells.Item(8, 10).Value = 0.3326292332200106
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 3.998462333333334
$ws.Cells.Item(8, 14).Value = 11.995387
$ws.Cells.Item(8, 15).Value = 0.07642924626252497
$ws.Cells.Item(8, 16).Value = 0.07778320264879995
$ws.Cells.Item(8, 17).Value = 7.168751152799667
$ws.Cells.Item(8, 18).Value = 64.518760375197
$ws.Cells.Item(8, 19).Value = 0.02542260157988704
$ws.Cells.Item(8, 20).Value = 0.02587296705446703

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.792877
$ws.Cells.Item(9, 8).Value = 5.378630999999999
$ws.Cells.Item(9, 9).Value = 0.3326292332200106
$ws.Cells.Item(9, 10).Value = 0.3326292332200106
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 21.74728766666667
$ws.Cells.Item(9, 14).Value = 65.241863
$ws.Cells.Item(9, 15).Value = 0.4156920000874432
$ws.Cells.Item(9, 16).Value = 0.4230560507063459
$ws.Cells.Item(9, 17).Value = 38.99021186995033
$ws.Cells.Item(9, 18).Value = 350.9119068295529
$ws.Cells.Item(9, 19).Value = 0.1382713112447788
$ws.Cells.Item(9, 20).Value = 0.1407208097555378

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.792877
$ws.Cells.Item(10, 8).Value = 5.378630999999999
$ws.Cells.Item(10, 9).Value = 0.3326292332200106
$ws.Cells.Item(10, 10).Value = 0.3326292332200106
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 20.37655766666667
$ws.Cells.Item(10, 14).Value = 61.129673
$ws.Cells.Item(10, 15).Value = 0.3894909627896643
$ws.Cells.Item(10, 16).Value = 0.3963908578200832
$ws.Cells.Item(10, 17).Value = 36.53266157974033
$ws.Cells.Item(10, 18).Value = 328.7939542176629
$ws.Cells.Item(10, 19).Value = 0.1295560802988497
$ws.Cells.Item(10, 20).Value = 0.1318511870921165

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.792877
$ws.Cells.Item(11, 8).Value = 5.378630999999999
$ws.Cells.Item(11, 9).Value = 0.3326292332200106
$ws.Cells.Item(11, 10).Value = 0.3326292332200106
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 2.731955
$ws.Cells.Item(11, 14).Value = 5.46391
$ws.Cells.Item(11, 15).Value = 0.0522203897564463
$ws.Cells.Item(11, 16).Value = 0.03543032157151783
$ws.Cells.Item(11, 17).Value = 4.898059284535
$ws.Cells.Item(11, 18).Value = 29.38835570721
$ws.Cells.Item(11, 19).Value = 0.01737002820313683
$ws.Cells.Item(11, 20).Value = 0.01178516069707238

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 3.311812
$ws.Cells.Item(12, 8).Value = 9.935435999999999
$ws.Cells.Item(12, 9).Value = 0.6144345017136311
$ws.Cells.Item(12, 10).Value = 0.6144345017136311
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 3.461605
$ws.Cells.Item(12, 14).Value = 10.384815
$ws.Cells.Item(12, 15).Value = 0.06616740110392129
$ws.Cells.Item(12, 16).Value = 0.06733956725325306
$ws.Cells.Item(12, 17).Value = 11.46418497826
$ws.Cells.Item(12, 18).Value = 103.17766480434
$ws.Cells.Item(12, 19).Value = 0.04065553412697384
$ws.Cells.Item(12, 20).Value = 0.0413757534508641

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 3.311812
$ws.Cells.Item(13, 8).Value = 9.935435999999999
$ws.Cells.Item(13, 9).Value = 0.6144345017136311
$ws.Cells.Item(13, 10).Value = 0.6144345017136311
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 3.998462333333334
$ws.Cells.Item(13, 14).Value = 11.995387
$ws.Cells.Item(13, 15).Value = 0.07642924626252497
$ws.Cells.Item(13, 16).Value = 0.07778320264879995
$ws.Cells.Item(13, 17).Value = 13.24215553708133
$ws.Cells.Item(13, 18).Value = 119.179399833732
$ws.Cells.Item(13, 19).Value = 0.04696076584366293
$ws.Cells.Item(13, 20).Value = 0.04779268336120579

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 3.311812
$ws.Cells.Item(14, 8).Value = 9.935435999999999
$ws.Cells.Item(14, 9).Value = 0.6144345017136311
$ws.Cells.Item(14, 10).Value = 0.6144345017136311
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 21.74728766666667
$ws.Cells.Item(14, 14).Value = 65.241863
$ws.Cells.Item(14, 15).Value = 0.4156920000874432
$ws.Cells.Item(14, 16).Value = 0.4230560507063459
$ws.Cells.Item(14, 17).Value = 72.02292826191865
$ws.Cells.Item(14, 18).Value = 648.2063543572679
$ws.Cells.Item(14, 19).Value = 0.2554155069400709
$ws.Cells.Item(14, 20).Value = 0.2599402337126903

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 3.311812
$ws.Cells.Item(15, 8).Value = 9.935435999999999
$ws.Cells.Item(15, 9).Value = 0.6144345017136311
$ws.Cells.Item(15, 10).Value = 0.6144345017136311
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 20.37655766666667
$ws.Cells.Item(15, 14).Value = 61.129673
$ws.Cells.Item(15, 15).Value = 0.3894909627896643
$ws.Cells.Item(15, 16).Value = 0.3963908578200832
$ws.Cells.Item(15, 17).Value = 67.48332819915866
$ws.Cells.Item(15, 18).Value = 607.3499537924279
$ws.Cells.Item(15, 19).Value = 0.2393166856436298
$ws.Cells.Item(15, 20).Value = 0.2435562192085216

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 3.311812
$ws.Cells.Item(16, 8).Value = 9.935435999999999
$ws.Cells.Item(16, 9).Value = 0.6144345017136311
$ws.Cells.Item(16, 10).Value = 0.6144345017136311
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 2.731955
$ws.Cells.Item(16, 14).Value = 5.46391
$ws.Cells.Item(16, 15).Value = 0.0522203897564463
$ws.Cells.Item(16, 16).Value = 0.03543032157151783
$ws.Cells.Item(16, 17).Value = 9.04772135246
$ws.Cells.Item(16, 18).Value = 54.28632811476
$ws.Cells.Item(16, 19).Value = 0.03208600915929369
$ws.Cells.Item(16, 20).Value = 0.02176961198034928
